# Add contact info for "Police Executive Research Forum" row in the
# Potential Partners table: Contact Name, Email, Phone.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(20)

# --- Contact Name cell (two runs: name + parenthetical note) ---
$nameCell = $row.Cells.Item(3)
$nameCell.Range.Text = "Nathan Ballard"

$nameParaRange = $d.Tables.Item(1).Rows.Item(20).Cells.Item(3).Range.Paragraphs.Item(1).Range.Duplicate
$nameParaRange.MoveEnd(1, -1)
$nameParaRange.Collapse(0)
$nameParaRange.Text = " (actually interviewing for a research position, but may refer to proper contact)"

# Format run 1 ("Nathan Ballard") with complex-script Arial font
$find1 = $d.Tables.Item(1).Rows.Item(20).Cells.Item(3).Range
$find1.Find.ClearFormatting()
$find1.Find.Text = "Nathan Ballard"
$find1.Find.Replacement.ClearFormatting()
$find1.Find.Replacement.Font.NameBi = "Arial"
$find1.Find.Replacement.Text = "Nathan Ballard"
$find1.Find.Execute($find1.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Find.Replacement.Text, 2)

# Format run 2 (" (actually interviewing ...)") with complex-script Arial font
$find2 = $d.Tables.Item(1).Rows.Item(20).Cells.Item(3).Range
$find2.Find.ClearFormatting()
$find2.Find.Text = " (actually interviewing for a research position, but may refer to proper contact)"
$find2.Find.Replacement.ClearFormatting()
$find2.Find.Replacement.Font.NameBi = "Arial"
$find2.Find.Replacement.Text = " (actually interviewing for a research position, but may refer to proper contact)"
$find2.Find.Execute($find2.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Find.Replacement.Text, 2)

# Set complex-script size (szCs=24 half-points => 12pt) for the whole cell
$d.Tables.Item(1).Rows.Item(20).Cells.Item(3).Range.Font.SizeBi = 12

# --- Email cell ---
$emailCell = $d.Tables.Item(1).Rows.Item(20).Cells.Item(4)
$emailCell.Range.Text = "nballard@policeforum.org"

$findEmail = $d.Tables.Item(1).Rows.Item(20).Cells.Item(4).Range
$findEmail.Find.ClearFormatting()
$findEmail.Find.Text = "nballard@policeforum.org"
$findEmail.Find.Replacement.ClearFormatting()
$findEmail.Find.Replacement.Font.NameBi = "Arial"
$findEmail.Find.Replacement.Text = "nballard@policeforum.org"
$findEmail.Find.Execute($findEmail.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $findEmail.Find.Replacement.Text, 2)

$d.Tables.Item(1).Rows.Item(20).Cells.Item(4).Range.Font.SizeBi = 12

# --- Phone cell ---
$phoneCell = $d.Tables.Item(1).Rows.Item(20).Cells.Item(5)
$phoneCell.Range.Text = "202-466-7820"

$findPhone = $d.Tables.Item(1).Rows.Item(20).Cells.Item(5).Range
$findPhone.Find.ClearFormatting()
$findPhone.Find.Text = "202-466-7820"
$findPhone.Find.Replacement.ClearFormatting()
$findPhone.Find.Replacement.Font.NameBi = "Arial"
$findPhone.Find.Replacement.Text = "202-466-7820"
$findPhone.Find.Execute($findPhone.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $findPhone.Find.Replacement.Text, 2)

$d.Tables.Item(1).Rows.Item(20).Cells.Item(5).Range.Font.SizeBi = 12
